$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9, column C ("MenuUIForm.Tips" row): update the translated value and
# give it the Chinese "微软雅黑" font, matching the localization update.
$cell = $ws.Range("C9")
$cell.Value = "空格键开始"
$cell.Font.Name = "微软雅黑"
$cell.Font.Size = 11
